$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price column D, Volume(1h) column E)
# D-column values that look numeric must be forced to Text so formats
# like trailing zeros / thousand-dot separators survive the COM write.

$ws.Range("D2").Value = '30.594.25'
$ws.Range("E2").Value = '  +1.56%  '
$ws.Range("D3").Value = '1.922.47'
$ws.Range("E3").Value = '  +3.83%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.59'
$ws.Range("E5").Value = '  +5.01%  '
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4729'
$ws.Range("E7").Value = '  +1.75%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2902'
$ws.Range("E8").Value = '  +4.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06791'
$ws.Range("E9").Value = '  +6.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '105.17'
$ws.Range("E10").Value = '  +8.63%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '18.38'
$ws.Range("E11").Value = '  +1.39%  '
$ws.Range("D12").Value = '1.911.64'
$ws.Range("E12").Value = '  +3.75%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07707'
$ws.Range("E13").Value = '  +2.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.292'
$ws.Range("E14").Value = '  +6.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6735'
$ws.Range("E15").Value = '  +7.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '288.00'
$ws.Range("E16").Value = '  -2.25%  '
$ws.Range("D17").Value = '30.611.26'
$ws.Range("E17").Value = '  +1.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007617'
$ws.Range("E18").Value = '  +3.08%  '
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.93'
$ws.Range("E20").Value = '  +1.91%  '
$ws.Range("D21").Value = '2.159.74'
$ws.Range("E21").Value = '  +3.81%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.441'
$ws.Range("E22").Value = '  +8.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9998'
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.330'
$ws.Range("E24").Value = '  +3.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.404'
$ws.Range("E25").Value = '  +3.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.06'
$ws.Range("E26").Value = '  +1.93%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.93'
$ws.Range("E27").Value = '  +8.73%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.122'
$ws.Range("E28").Value = '  +9.64%  '
$ws.Range("E29").Value = '  +0.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.369'
$ws.Range("E30").Value = '  +3.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.180'
$ws.Range("E31").Value = '  +4.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.138'
$ws.Range("E32").Value = '  +8.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05045'
$ws.Range("E33").Value = '  +2.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7421'
$ws.Range("E34").Value = '  +2.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.155'
$ws.Range("E35").Value = '  +3.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02078'
$ws.Range("E36").Value = '  +9.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.746'
$ws.Range("E37").Value = '  +0.58%  '
$ws.Range("E38").Value = '  +1.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.057'
$ws.Range("E39").Value = '  +4.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '111.06'
$ws.Range("E40").Value = '  +6.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8821'
$ws.Range("E41").Value = '  +2.48%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4358'
$ws.Range("E42").Value = '  +7.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.901'
$ws.Range("E43").Value = '  +4.48%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '67.19'
$ws.Range("E45").Value = '  +2.82%  '
$ws.Range("E46").Value = '  +2.87%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.244'
$ws.Range("E47").Value = '  +3.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '47.97'
$ws.Range("E48").Value = '  +17.40%  '
$ws.Range("E49").Value = '  +3.71%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.92'
$ws.Range("E50").Value = '  +2.79%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4039'
$ws.Range("E51").Value = '  +8.46%  '
